$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 3.35
$ws.Range("I2").Value = 3.45
$ws.Range("K2").Value = 3.3
$ws.Range("S2").Value = 4.7
$ws.Range("U2").Value = 1.94
$ws.Range("AA2").Value = 65
$ws.Range("AD2").Value = 14.5
$ws.Range("AE2").Value = 46
$ws.Range("AF2").Value = 15.5
$ws.Range("AO2").Value = 55

# Row 3
$ws.Range("H3").Value = 9.6
$ws.Range("R3").Value = 1.32

# Row 4
$ws.Range("J4").Value = 4.4
$ws.Range("K4").Value = 4.5
$ws.Range("Q4").Value = 1.82
$ws.Range("U4").Value = 2.1
$ws.Range("AH4").Value = 19.5
